$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (group-size values) B1:E1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# CON measurements, row 2, B:E
$ws.Range("B2").Value = 172.67528433332257
$ws.Range("C2").Value = 223.01334209369301
$ws.Range("D2").Value = 174.58653062100834
$ws.Range("E2").Value = 223.133826288135

# STR measurements, row 3, B:E
$ws.Range("B3").Value = 171.65790081672372
$ws.Range("C3").Value = 227.08603722823111
$ws.Range("D3").Value = 180.64958327106208
$ws.Range("E3").Value = 217.89594454113043

# Match the author's updated selection extent (B1:E3 instead of the full B1:AY3)
$ws.Range("B1:E3").Select()
